$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Section 1 (CustomerMapping Class) grading points, copying the
# "Points for grading" column (D) values into the "Grading comments" (E) column
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Fill in Section 2 (Customer Class) grading points similarly
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the active selection to E15, matching where grading left off
$ws.Range("E15").Select()

$wb.Save()
